$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017年")

# Extend the formatting of the data rows (copy formats from the row above,
# matching how the table had been filled in for every prior row) down into
# the two new rows, 22 and 23.
$ws.Range("B21:G21").Copy() | Out-Null
$ws.Range("B22:G23").PasteSpecial(-4122) | Out-Null

# Row 22: 支出, 生活费 300 on 2017-11-21, note 生活费(11/21-11/30)
$ws.Cells.Item(22, 3).Value = "支出"
$ws.Cells.Item(22, 4).Value = 300
$ws.Cells.Item(22, 5).Value = 43060
$ws.Cells.Item(22, 6).Value = "生活费"
$ws.Cells.Item(22, 7).Value = "生活费(11/21-11/30)"

# Row 23: 支出, 班费(其他) 100 on 2017-11-21, note 班费
$ws.Cells.Item(23, 3).Value = "支出"
$ws.Cells.Item(23, 4).Value = 100
$ws.Cells.Item(23, 5).Value = 43060
$ws.Cells.Item(23, 6).Value = "其他"
$ws.Cells.Item(23, 7).Value = "班费"

$ws.Range("G24").Select()
